# "add config for chapter"
#
# The Chapters sheet used to carry a second, unrelated config block in
# columns F:H (priceMultiplier / priceBase / buildingList, with sample
# values + the "affect the building price" / "reward money amount"
# notes). That block is removed here, and the chapter IDs in column B
# are renumbered from the placeholder 1001..1006 scheme down to the
# plain 1..6 scheme used everywhere else (Buildings sheet, etc).

$wb = $excel.ActiveWorkbook
$chapters = $wb.Worksheets.Item("Chapters")
$buildings = $wb.Worksheets.Item("Buildings")

# Drop the extra price/building-list config columns (F:H) from the
# header/meta rows 1-3, keeping row 1/2's styling on F/G/H but clearing
# the text, and dropping the now-empty-styled F3/G3 cells outright.
$chapters.Range("F1").ClearContents()
$chapters.Range("G1").ClearContents()
$chapters.Range("H1").ClearContents()
$chapters.Range("F2").ClearContents()
$chapters.Range("G2").ClearContents()
$chapters.Range("H2").ClearContents()
$chapters.Range("F3").ClearContents()
$chapters.Range("G3").ClearContents()

# ... and the matching per-chapter data in F4:H9.
$chapters.Range("F4:H9").ClearContents()

# Renumber the chapter IDs (B4:B9) from 1001..1006 to 1..6.
$chapters.Range("B4").Value = 1
$chapters.Range("B5").Value = 2
$chapters.Range("B6").Value = 3
$chapters.Range("B7").Value = 4
$chapters.Range("B8").Value = 5
$chapters.Range("B9").Value = 6

# Leave the workbook with the Chapters sheet active/selected, matching
# the edited file's last-saved view state.
$buildings.Activate()
$buildings.Range("C6:E20").Select()
$chapters.Activate()
$chapters.Range("E7").Select()
